$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.125.24'
$ws.Range("E2").Value = '  -0.18%  '

$ws.Range("D3").Value = '3.424.55'
$ws.Range("E3").Value = '  +1.42%  '

$ws.Range("E4").Value = '  -0.10%  '

$ws.Range("D5").Value = "'547.94"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.35%  '

$ws.Range("D6").Value = "'178.02"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.62%  '

$ws.Range("D7").Value = "'0.635"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +5.71%  '

$ws.Range("E8").Value = '  +0.00%  '

$ws.Range("D9").Value = "'0.623"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.27%  '

$ws.Range("D10").Value = "'0.149"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +6.82%  '

$ws.Range("D11").Value = "'53.31"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.90%  '

$ws.Range("E12").Value = '  +1.37%  '

$ws.Range("D13").Value = "'9.12"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.64%  '

$ws.Range("D14").Value = '3.971.50'
$ws.Range("E14").Value = '  +1.22%  '

$ws.Range("E15").Value = '  +0.80%  '

$ws.Range("D16").Value = '3.424.37'
$ws.Range("E16").Value = '  +1.30%  '

$ws.Range("D17").Value = "'18.19"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.81%  '

$ws.Range("D18").Value = '65.236.37'
$ws.Range("E18").Value = '  -0.40%  '

$ws.Range("D19").Value = "'11.76"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.29%  '

$ws.Range("E20").Value = '  -0.49%  '

$ws.Range("D21").Value = "'412.59"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +6.49%  '

$ws.Range("E22").Value = '  +5.59%  '

$ws.Range("B23").Value = 'Litecoin'
$ws.Range("C23").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D23").Value = "'84.79"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.38%  '

$ws.Range("B24").Value = 'Toncoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D24").Value = "'4.10"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.98%  '

$ws.Range("D25").Value = "'10.71"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -9.91%  '

$ws.Range("E26").Value = '  +0.72%  '

$ws.Range("E27").Value = '  +5.44%  '

$ws.Range("D28").Value = "'6.03"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.46%  '

$ws.Range("E29").Value = '  +4.85%  '

$ws.Range("D30").Value = "'29.62"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.77%  '

$ws.Range("D31").Value = "'6.48"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -4.81%  '

$ws.Range("D32").Value = "'608.55"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -7.67%  '

$ws.Range("D33").Value = "'11.59"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.42%  '

$ws.Range("E34").Value = '  +0.20%  '

$ws.Range("D35").Value = "'58.85"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.54%  '

$ws.Range("B36").Value = 'Dai'
$ws.Range("C36").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D36").Value = "'0.999"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.11%  '

$ws.Range("B37").Value = 'Kaspa'
$ws.Range("C37").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D37").Value = "'0.146"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +13.25%  '

$ws.Range("E38").Value = '  -1.11%  '

$ws.Range("E39").Value = '  -1.44%  '

$ws.Range("E40").Value = '  -4.50%  '

$ws.Range("D41").Value = '3.175.15'
$ws.Range("E41").Value = '  +5.43%  '

$ws.Range("E42").Value = '  +1.81%  '

$ws.Range("D43").Value = "'1.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.25%  '

$ws.Range("D44").Value = "'2.52"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -8.88%  '

$ws.Range("B45").Value = 'ApeXProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D45").Value = "'3.27"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.90%  '

$ws.Range("B46").Value = 'ThetaToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D46").Value = "'2.77"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.89%  '

$ws.Range("D47").Value = "'0.0407"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.46%  '

$ws.Range("E48").Value = '  +0.38%  '

$ws.Range("E49").Value = '  +2.98%  '

$ws.Range("D50").Value = "'137.84"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.10%  '

$ws.Range("D51").Value = "'8.34"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.85%  '
